$d = $word.ActiveDocument

$replacements = @(
    @("628×6=3768", "113×2=226"),
    @("556×5=2780", "875×9=7875"),
    @("117×2=234", "399×4=1596"),
    @("775×6=4650", "459×3=1377"),
    @("679×3=2037", "432×3=1296"),
    @("690×9=6210", "146×2=292"),
    @("690×8=5520", "972×5=4860"),
    @("982×9=8838", "941×3=2823"),
    @("997×7=6979", "372×3=1116"),
    @("444×8=3552", "683×3=2049"),
    @("420×7=2940", "820×6=4920"),
    @("509×5=2545", "703×7=4921"),
    @("803×9=7227", "997×8=7976"),
    @("427×4=1708", "408×2=816"),
    @("819×2=1638", "958×4=3832"),
    @("670×9=6030", "423×9=3807"),
    @("410×8=3280", "209×8=1672"),
    @("509×6=3054", "422×6=2532"),
    @("928×9=8352", "991×5=4955"),
    @("601×3=1803", "394×3=1182"),
    @("401×8=3208", "986×9=8874"),
    @("475×5=2375", "348×3=1044"),
    @("404×8=3232", "403×7=2821"),
    @("135×8=1080", "300×2=600"),
    @("814×9=7326", "732×7=5124")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
